$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 entirely (it no longer belongs in the sheet).
$ws.Rows.Item(3).Delete()

# Capture the existing Player 5..8 data currently living in row 2
# (columns AJ..BD, i.e. column indices 36..56) before it gets overwritten,
# since it needs to slide right by 7 columns to make room for a brand new
# "Player 5" block.
$oldVals = @()
for ($c = 36; $c -le 56; $c++) {
    $oldVals += ,$ws.Cells.Item(2, $c).Value()
}

# Move that captured data into its new home: columns AQ..BK (indices 43..63).
for ($i = 0; $i -lt $oldVals.Count; $i++) {
    $ws.Cells.Item(2, 43 + $i).Value = $oldVals[$i]
}

# Fill in the freshly inserted "Player 5" columns (AJ..AP, indices 36..42)
# with the new warning-message related match data.
$newPlayer5 = @("Lumi#2340", "N-ZAP '85", "6", "2", "4", "1", "463")
for ($i = 0; $i -lt $newPlayer5.Count; $i++) {
    $ws.Cells.Item(2, 36 + $i).Value = $newPlayer5[$i]
}
